$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Automatic report generation"
$ws.Range("C9").Value = "Report"
$ws.Range("D9").Value = "Open"

$ws.Range("B9").Select()
